$d = $word.ActiveDocument

# 1. Insert " sourced from the streaming servers" right after
#    "...7.5Gbps of unicast traffic" (and before the following period).
$d.Content.Find.Execute(
    "500 households can amount to 7.5Gbps of unicast traffic.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "500 households can amount to 7.5Gbps of unicast traffic sourced from the streaming servers.",
    2) | Out-Null

# 2. "This amounts to 45Gbps from the cache cluster" -> "... the streaming server cluster"
$d.Content.Find.Execute(
    "This amounts to 45Gbps from the cache cluster and 37.5Gbps",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This amounts to 45Gbps from the streaming server cluster and 37.5Gbps",
    2) | Out-Null

# 3. Word keeps the "_GoBack" bookmark pinned to the spot of the most recent
#    edit. Relocate it (collapsed) to right after "...streaming server",
#    i.e. right before " cluster and 37.5Gbps..." -- matching where the
#    last text change above landed.
$locate = $d.Content
$locate.Find.Execute(
    "This amounts to 45Gbps from the streaming server",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$goBackPos = $locate.End

$existing = $d.Bookmarks.Item("_GoBack")
$existing.Delete()

$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
